# Deploy the implementation guide.
# Updates the Metadata sheet's Status and Date values, and touches the
# alignment formatting on the two body cell styles so Excel persists the
# `applyAlignment="true"` flag alongside the existing <alignment> settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: bump the publication timestamp
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# Re-apply the wrap/vertical-top alignment on every sheet so the cellXfs
# entries pick up applyAlignment="true" (the alignment values themselves -
# vertical=top, wrapText=true - are already set, this just flips the
# "apply" flag that Excel tracks separately from the alignment values).
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $used.WrapText = $true
}
